$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values: force text format so purely numeric-looking
# strings (e.g. "589.89") stay text cells like the source data, instead of
# being auto-coerced to Excel numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.645.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.127.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.116.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.648.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.580.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.124.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.735'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0861'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '447.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0372'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.920.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.58'
$ws.Range("D51").Style = "Normal"

# Volume(1h) (column E) values: plain text assignment (never parses as a number
# due to surrounding spaces and percent sign).
$ws.Range("E2").Value = '  +3.12%  '
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("E6").Value = '  +2.93%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("E10").Value = '  +17.09%  '
$ws.Range("E11").Value = '  +3.80%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  +6.63%  '
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("E18").Value = '  +3.17%  '
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("E20").Value = '  +3.58%  '
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("E23").Value = '  +1.47%  '
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  +8.31%  '
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  +8.63%  '
$ws.Range("E35").Value = '  +10.47%  '
$ws.Range("E36").Value = '  +2.01%  '
$ws.Range("E37").Value = '  +13.76%  '
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("E40").Value = '  +5.76%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("E43").Value = '  +5.29%  '
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("E46").Value = '  +3.40%  '
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("E51").Value = '  +3.07%  '
